$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo: "éléménets" -> "éléments" in the cart-display row's "how to verify" cell (E12)
$ws.Range("E12").Value = "les éléments s'affichent correctement et completement : l'image, son texte alt, le nom du kanap, sa couleur, son prix, sa quantité désirée, son prix total, le nombre d'articles total"

# Fix typo: "reload" -> "recharge" in the deleteItem() row's "expected result" cell (D10)
$ws.Range("D10").Value = "supprime une entrée du panier de html et du localStorage, et recharge la page"

# Widen column E to fit the relocated text
$ws.Columns("E").ColumnWidth = 105.83072916666667

# Move the active selection to D10 (matches the saved view state in the workbook)
$ws.Range("D10").Select()
